$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Статистика по годам")

function Swap-Row($rowA, $rowB) {
    for ($col = 1; $col -le 5; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)
        $valA = $cellA.Value2
        $valB = $cellB.Value2
        $cellA.Value2 = $valB
        $cellB.Value2 = $valA
    }
}

# Swap row 3 and row 4 (columns A:E)
Swap-Row 3 4

# Swap row 6 and row 7 (columns A:E)
Swap-Row 6 7
